# The workbook tracks daily price observations for "Repollo" (cabbage) at
# Feria Lagunitas de Puerto Montt. Two new daily observations need to be
# inserted right before the current row 159, pushing the existing rows
# 159-246 down to 161-248 (dimension grows from A1:R246 to A1:R248).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 159 (Excel copies the
# formatting, e.g. the date number format in column D, from the row above).
$ws.Rows.Item(159).Insert()
$ws.Rows.Item(159).Insert()

# New row 159: Copenhague / Primera, date serial 44460 (2021-09-21)
$ws.Cells.Item(159, 1).Value = 4
$ws.Cells.Item(159, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(159, 3).Value = "Los Lagos"
$ws.Cells.Item(159, 4).Value = 44460
$ws.Cells.Item(159, 5).Value = 10
$ws.Cells.Item(159, 6).Value = 100112006
$ws.Cells.Item(159, 7).Value = "Repollo"
$ws.Cells.Item(159, 8).Value = "Copenhague"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 400
$ws.Cells.Item(159, 11).Value = 1500
$ws.Cells.Item(159, 12).Value = 1500
$ws.Cells.Item(159, 13).Value = 1500
$ws.Cells.Item(159, 14).Value = "$/unidad"
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 1500
$ws.Cells.Item(159, 17).Value = 1
$ws.Cells.Item(159, 18).Value = "Hortaliza"

# New row 160: Crespo record / Segunda, same date serial 44460
$ws.Cells.Item(160, 1).Value = 4
$ws.Cells.Item(160, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(160, 3).Value = "Los Lagos"
$ws.Cells.Item(160, 4).Value = 44460
$ws.Cells.Item(160, 5).Value = 10
$ws.Cells.Item(160, 6).Value = 100112006
$ws.Cells.Item(160, 7).Value = "Repollo"
$ws.Cells.Item(160, 8).Value = "Crespo record"
$ws.Cells.Item(160, 9).Value = "Segunda"
$ws.Cells.Item(160, 10).Value = 500
$ws.Cells.Item(160, 11).Value = 1000
$ws.Cells.Item(160, 12).Value = 1000
$ws.Cells.Item(160, 13).Value = 1000
$ws.Cells.Item(160, 14).Value = "$/unidad"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 1000
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"
